$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at rows 3-5 (pushes existing rows 3.. down to 6..)
$ws.Rows("3:5").Insert()

# Row 3: Date Time
$ws.Range("A3").Value = "Date Time"
$ws.Range("B3").Value = 42370.5
$ws.Range("B3").NumberFormat = "m/d/yy h:mm"

# Row 4: Time
$ws.Range("A4").Value = "Time"
$ws.Range("B4").Value = 0.5
$ws.Range("B4").NumberFormat = "h:mm AM/PM"

# Row 5: Midnight
$ws.Range("A5").Value = "Midnight"
$ws.Range("B5").Value = 0
$ws.Range("B5").NumberFormat = "h:mm AM/PM"

# Resize column B to fit its new (wider) contents
$ws.Columns("B:B").ColumnWidth = 10.83

# Update the selected cell to match the new edit location
$ws.Range("B6").Select()
